# Daily attendance processing - reorder "Recorded By" names in column G
# so that "System"/"system" is listed first among the comma-separated
# recorder names (matches the canonical recorder-name ordering).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -match ",") {
        $parts = $val -split ", "

        if ($parts.Count -gt 1) {
            $first = $parts[0]

            if ($first.ToLower() -ne "system") {
                $hasSystem = $false
                foreach ($p in $parts) {
                    if ($p.ToLower() -eq "system") {
                        $hasSystem = $true
                    }
                }

                if ($hasSystem) {
                    $reversed = $parts[($parts.Count - 1)..0]
                    $cell.Value = ($reversed -join ", ")
                }
            }
        }
    }
}
